$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by 3 rows (39 -> 42), copying formatting from the last
# existing data row (38) so the new rows 40-42 match the established look.
$ws.Range("A38:E38").Copy($ws.Range("A40:E40"))
$ws.Range("A38:E38").Copy($ws.Range("A41:E41"))
$ws.Range("A38:E38").Copy($ws.Range("A42:E42"))

# Order numbers for the new rows
$ws.Range("A40").Value = 39
$ws.Range("A41").Value = 40
$ws.Range("A42").Value = 41

# Rewrite Field Name / Field Category / Field Type / Description for rows 15-42
$ws.Range("B15").Value = "DecisionExplicit"
$ws.Range("C15").Value = "Decision Analysis"
$ws.Range("D15").Value = "Yes/No"
$ws.Range("E15").Value = "Was the decision of the process documented?"

$ws.Range("B16").Value = "OptimalAltExplicit"
$ws.Range("C16").Value = "Decision Analysis"
$ws.Range("D16").Value = "Yes/No"
$ws.Range("E16").Value = "Was the best management procedure, aka optimal alternative, documented?"

$ws.Range("B17").Value = "RoleSpecification"
$ws.Range("C17").Value = "Decision Process"
$ws.Range("D17").Value = "Yes/No"
$ws.Range("E17").Value = "Were roles assigned and documented?"

$ws.Range("B18").Value = "OpenMeetings"
$ws.Range("C18").Value = "Decision Process"
$ws.Range("D18").Value = "Yes/No"
$ws.Range("E18").Value = "Were open meetings held?"

$ws.Range("B19").Value = "ResultsAdopted"
$ws.Range("C19").Value = "Decision Process"
$ws.Range("D19").Value = "Yes/No"
$ws.Range("E19").Value = "Did the MSE influence subsequent management?"

$ws.Range("B20").Value = "ProblemDefinition"
$ws.Range("C20").Value = "Decision Analysis"
$ws.Range("D20").Value = "Description"
$ws.Range("E20").Value = "A problem definition taken from the documentation (Reader interpretation)"

$ws.Range("B21").Value = "ObjElicitationMethod"
$ws.Range("C21").Value = "Decision Analysis"
$ws.Range("D21").Value = "Description"
$ws.Range("E21").Value = "If documented, how were objectives elicited?"

$ws.Range("B22").Value = "TradeOffMethod_Exp"
$ws.Range("C22").Value = "Decision Analysis"
$ws.Range("D22").Value = "Description"
$ws.Range("E22").Value = "If explicitly documented, what form of tradeoff analysis occured?"

$ws.Range("B23").Value = "TradeOffMethod_Sub"
$ws.Range("C23").Value = "Decision Analysis"
$ws.Range("D23").Value = "Description"
$ws.Range("E23").Value = "If not explicitly documented, what form of tradeoff analysis seems to have occured?"

$ws.Range("B24").Value = "Decision"
$ws.Range("C24").Value = "Result"
$ws.Range("D24").Value = "Description"
$ws.Range("E24").Value = "If documented, the management procedure that was selected for implementation"

$ws.Range("B25").Value = "Leader"
$ws.Range("C25").Value = "Decision Process"
$ws.Range("D25").Value = "List"
$ws.Range("E25").Value = "What organization initiated and directed the MSE?"

$ws.Range("B26").Value = "Participants"
$ws.Range("C26").Value = "Decision Process"
$ws.Range("D26").Value = "List"
$ws.Range("E26").Value = "Who participated in the MSE process?"

$ws.Range("B27").Value = "ObjElicitationSource_Exp"
$ws.Range("C27").Value = "Decision Process"
$ws.Range("D27").Value = "List"
$ws.Range("E27").Value = "If explicitly documented, the groups from which objectives were elicited"

$ws.Range("B28").Value = "ObjElicitationSource_Sub"
$ws.Range("C28").Value = "Decision Process"
$ws.Range("D28").Value = "List"
$ws.Range("E28").Value = "If not explicitly documented, the groups from which objectives were seemingly elicited"

$ws.Range("B29").Value = "ProcedureElicitation_Exp"
$ws.Range("C29").Value = "Decision Process"
$ws.Range("D29").Value = "List"
$ws.Range("E29").Value = "If explicitly documented, the groups from which alternative management procedures were elicited"

$ws.Range("B30").Value = "ProcedureElicitation_Sub"
$ws.Range("C30").Value = "Decision Process"
$ws.Range("D30").Value = "List"
$ws.Range("E30").Value = "If not explicitly documented, the groups from which alternative management procedures were seemingly elicited"

$ws.Range("B31").Value = "ConsequencePrediction"
$ws.Range("C31").Value = "Decision Process"
$ws.Range("D31").Value = "Description"
$ws.Range("E31").Value = "How were consequences predicted?"

$ws.Range("B32").Value = "FullCitation"
$ws.Range("C32").Value = "Documentation"
$ws.Range("D32").Value = "Description"
$ws.Range("E32").Value = "The full citation for the study"

$ws.Range("B33").Value = "Comments"
$ws.Range("C33").Value = "Comments"
$ws.Range("D33").Value = "Description"
$ws.Range("E33").Value = "Additional notes and comments about the study"

$ws.Range("B34").Value = "ObjName"
$ws.Range("C34").Value = "Objectives"
$ws.Range("D34").Value = "Description"
$ws.Range("E34").Value = "Text description of the objective"

$ws.Range("B35").Value = "ObjCategory"
$ws.Range("C35").Value = "Objectives"
$ws.Range("D35").Value = "List"
$ws.Range("E35").Value = "The objective category (conservation, yield, economic, social)"

$ws.Range("B36").Value = "ObjDescription"
$ws.Range("C36").Value = "Objectives"
$ws.Range("D36").Value = "Description"
$ws.Range("E36").Value = "Description of the objective"

$ws.Range("B37").Value = "ObjDirection"
$ws.Range("C37").Value = "Objectives"
$ws.Range("D37").Value = "List"
$ws.Range("E37").Value = "The desired state of the objective"

$ws.Range("B38").Value = "ObjType"
$ws.Range("C38").Value = "Objectives"
$ws.Range("D38").Value = "List"
$ws.Range("E38").Value = "The type of objective.  E.g., strategic, process, fundamental, or means"

$ws.Range("B39").Value = "ObjScale"
$ws.Range("C39").Value = "Objectives"
$ws.Range("D39").Value = "Description"
$ws.Range("E39").Value = "The scale on which the objective is measures (natural, proxy, or constructed)"

$ws.Range("B40").Value = "ObjMetric"
$ws.Range("C40").Value = "Objectives"
$ws.Range("D40").Value = "Description"
$ws.Range("E40").Value = "The units used to measure the objective"

$ws.Range("B41").Value = "MPManagementTool"
$ws.Range("C41").Value = "Alternatives"
$ws.Range("D41").Value = "List"
$ws.Range("E41").Value = "Management alternatives evaluated in the study"

$ws.Range("B42").Value = "MPAlternativesEvaluated"
$ws.Range("C42").Value = "Alternatives"
$ws.Range("D42").Value = "List"
$ws.Range("E42").Value = "Types of alternatives evaluated"

Write-Output "Workbook updated"